$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Old = 440938; New = 606323 },
    @{ Row = 3; Old = 175826; New = 503433 },
    @{ Row = 4; Old = 413334; New = 257987 },
    @{ Row = 5; Old = 837572; New = 799241 },
    @{ Row = 6; Old = 109417; New = 796294 },
    @{ Row = 7; Old = 483721; New = 184961 },
    @{ Row = 8; Old = 487282; New = 502539 },
    @{ Row = 9; Old = 929571; New = 929043 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $newId = $u.New
    $ws.Cells.Item($r, 4).Value = $newId
    $ws.Cells.Item($r, 5).Value = "https://93mtzf.deta.dev/ticket/$newId"
}
